$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from an existing "Easy" row (row 40) down into the two new rows,
# so the new rows reuse the existing style indices (fills/fonts) instead of Excel
# creating brand-new style records for them.
$ws.Range("A40:E40").Copy() | Out-Null
$ws.Range("A41:E41").PasteSpecial(-4122) | Out-Null
$ws.Range("A42:E42").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 41: 1071. Greatest Common Divisor of Strings
$ws.Range("A41").Value = "1071. Greatest Common Divisor of Strings"
$ws.Range("B41").Value = "Easy"
$ws.Range("C41").Value = "Strings"
$ws.Range("D41").Value = "First check if GCD is possible with str1 + str2 == str2 + str1. We use the Euclidean algorithm to recursively find the substring."
$ws.Hyperlinks.Add($ws.Range("E41"), "https://leetcode.com/problems/greatest-common-divisor-of-strings/solutions/3124997/super-easy-solution-fully-explained-c-python3-java/?envType=study-plan-v2&envId=leetcode-75") | Out-Null

# Row 42: 605. Can Place Flowers
$ws.Range("A42").Value = "605. Can Place Flowers"
$ws.Range("B42").Value = "Easy"
$ws.Range("C42").Value = "Arrays"
$ws.Range("D42").Value = "The Optimal solution is the Greedy solution. Greedily place a flower at every vacant spot from left to right. The Naïve solution is that you can do a for loop, checking each triplet if they are all 0's, and handle th edge cases (n=0, l=1, l=2, starting i and ending i. Increment a count."
$ws.Hyperlinks.Add($ws.Range("E42"), "https://leetcode.com/problems/can-place-flowers/solutions/103898/java-greedy-solution-o-flowerbed-beats-100/?envType=study-plan-v2&envId=leetcode-75 ") | Out-Null

# Hyperlinks.Add swaps in a slightly different (but visually identical) style
# variant for the link cells. Re-paste the original link-cell formatting so the
# two new cells collapse back onto the very same style index used by every
# other link cell in the sheet.
$ws.Range("E40").Copy() | Out-Null
$ws.Range("E41").PasteSpecial(-4122) | Out-Null
$ws.Range("E42").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# PasteSpecial(formats) does not touch cell values/hyperlinks, but just in case,
# make sure the display text of the link cells is still correct.
$ws.Range("E41").Value = "https://leetcode.com/problems/greatest-common-divisor-of-strings/solutions/3124997/super-easy-solution-fully-explained-c-python3-java/?envType=study-plan-v2&envId=leetcode-75"
$ws.Range("E42").Value = "https://leetcode.com/problems/can-place-flowers/solutions/103898/java-greedy-solution-o-flowerbed-beats-100/?envType=study-plan-v2&envId=leetcode-75 "

# Match the scrolled/selected view state recorded in the saved workbook.
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("E48").Select() | Out-Null
